$p = $ppt.ActivePresentation

# -------------------------------------------------------------------
# 1) Slide 5 - "Methodology" SmartArt diagram: the first node's text
#    was split across two runs ("Load dataset " + "from the source
#    into R"). Merge it back into a single logical string so the
#    diagram (and its cached drawing) shows one continuous run.
# -------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$diagramShape = $s5.Shapes.Item(2)
$smartArt = $diagramShape.SmartArt
$firstNode = $smartArt.Nodes.Item(1)
$firstNode.TextFrame.TextRange.Text = "Load dataset from the source into R"

# -------------------------------------------------------------------
# 2) Slide 8 - bullet list describing "Create tables, database".
#    a) "Manually" -> "Manually (3 tables)"
#    b) "Based on the ID" -> "Based on the ID, " followed by a new
#       run "code name (10 tables)" (same look/formatting).
# -------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$bodyShape = $s8.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

# a) "Manually" paragraph (paragraph 1 in this text box).
$manuallyPara = $bodyRange.Paragraphs(1, 1)
$manuallyFull = $manuallyPara.Characters(1, $manuallyPara.Length)
$manuallyFull.Text = "Manually (3 tables)"

# b) "Based on the ID" paragraph (paragraph 3 in this text box).
$bodyRange = $bodyShape.TextFrame.TextRange
$basedPara = $bodyRange.Paragraphs(3, 1)
$basedFull = $basedPara.Characters(1, $basedPara.Length)
$basedFull.Text = "Based on the ID, "
$basedPara.InsertAfter("code name (10 tables)")
